$wb = $excel.ActiveWorkbook

# Rename sheets: "Good" -> "Biomass" (formulas referencing these sheets are
# updated automatically by the host when the sheet is renamed).
$wsFlow = $wb.Worksheets.Item("Flow_data_Good")
$wsFlow.Name = "Flow_data_Biomass"

$wsTC = $wb.Worksheets.Item("TC_data_Good")
$wsTC.Name = "TC_data_Biomass"

# Make TC_data_Biomass (3rd sheet, index 2) the active / selected tab,
# and move its selection to R28.
$wsTC.Activate()
[void]$wsTC.Range("R28").Select()
